# Scraper update - Linea 141: new scrape pass at 01:30:59 adds a fresh
# top row to LP1912 ("15_ABASTO") plus a first-ever "215_ALUAR" stop that
# also spawns data on the LP1912-215 sheet. Sheet 6203-6173 stays empty,
# only its "last updated" stamp moves forward.

$wb = $excel.ActiveWorkbook
$newTime = "01:30:59"

# ----- Sheet 1: LP1912 -----
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: $newTime"
$ws1.Range("A3").Value = "Total filas: 2"

# Row 6 refreshed with the latest arrival for this pass
$ws1.Range("A6").Value = $newTime
$ws1.Range("B6").Value = "03:01"
$ws1.Range("C6").Value = "15_ABASTO"
$ws1.Range("D6").Value = 91
$ws1.Range("E6").Value = "LP1912"

# Row 7 newly appended
$ws1.Range("A7").Value = $newTime
$ws1.Range("B7").Value = "03:06"
$ws1.Range("C7").Value = "215_ALUAR"
$ws1.Range("D7").Value = 96
$ws1.Range("E7").Value = "LP1912"

# ----- Sheet 2: LP1912-215 -----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: $newTime"
$ws2.Range("A3").Value = "Total filas: 1"

# Bring in the header row formatting from sheet 1 (bold, bordered, centered)
$ws1.Range("A5:E5").Copy()
$ws2.Range("A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws2.Range("A5").Value = "Hora_Scrap"
$ws2.Range("B5").Value = "Hora_Llegada"
$ws2.Range("C5").Value = "Linea"
$ws2.Range("D5").Value = "Minutos"
$ws2.Range("E5").Value = "Parada"

$ws2.Range("A6").Value = $newTime
$ws2.Range("B6").Value = "03:06"
$ws2.Range("C6").Value = "215_ALUAR"
$ws2.Range("D6").Value = 96
$ws2.Range("E6").Value = "LP1912"

# ----- Sheet 3: 6203-6173 -----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: $newTime"
